$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1) "TextBox 28" (shape holding the "    ???? DB" label) -> "    SQl Server"
#    The leading 4 spaces keep their original 14pt size; the new "SQl Server"
#    text is 16pt/bold. The box auto-fits (spAutoFit), so its height must be
#    set *after* the text/font edits, otherwise autofit recomputes it from
#    the (taller) original text.
$tb = $s.Shapes.Item(6)
$tr = $tb.TextFrame.TextRange
$tr.Font.Size = 16
$tr.Font.Bold = $true
$tr.Text = "    SQl Server"
$tr.Characters(1, 4).Font.Size = 14
$tb.Height = 26.6578

# 2) "Rectangle 69" - shrink the big translucent overlay rectangle down to a
#    small marker box positioned over the "SQl Server" label.
$rect = $s.Shapes.Item(23)
$rect.Left = 162.5454331
$rect.Top = 449.7656
$rect.Width = 30.2902
$rect.Height = 46.85874

# 3) Remove the stray vertical "??" label ("TextBox 3", rotated textbox).
$s.Shapes.Item(28).Delete()
